# Quarterly indexing esoteric bug-fix operation
#
# Column A (rows 2..73) holds a date serial that should represent the
# "as-of" quarter reference date, but was mistakenly written as the 1st of
# the month.  The fix moves each such date forward to the 15th of the
# *following* month (i.e. EDATE(old,1) + 14 days), leaving every other
# value on the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -ne $null) {
        $shifted = $excel.WorksheetFunction.EDate($old, 1)
        $new = $shifted + 14
        $cell.Value2 = $new
    }
}
